$wb = $excel.ActiveWorkbook

# The existing last sheet ("Sheet17") is the template for the new check/test sheet.
$src = $wb.Worksheets.Item("Sheet17")

# Before moving away from Sheet17, widen its lingering selection from the single
# cell D5 to the sheet's whole used range (A1:J7), matching the post-edit state.
$src.Activate()
$src.Range("A1:J7").Select()

# Add the new worksheet right after Sheet17 -> becomes "Sheet18" / sheetId 18,
# and is placed at the end of the tab strip (becoming the new active tab).
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $src)

# Copy Sheet17's data (A1:J7) verbatim into the new sheet.
$src.Range("A1:J7").Copy()
$new.Range("A1").PasteSpecial()

# Leave the new sheet's selection where the author left it.
$new.Range("E27").Select()
